# "Generate Report for Handoff": a new handoff round was just generated for
# d3aa5e99-a3c6-4e8c-88fe-952729efaa67.md, so its "Latest Handoff
# Datetime" / "Latest HO Xliff Generate Date" timestamps (row 7 on each
# sheet) are refreshed to the new generation time.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-30 20:47:16"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-30 20:47:11"

# de-de sheet: column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-30 20:47:16"
